$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 33.17258048057556
$ws.Range("F3").Value = 32.05534434318542
$ws.Range("F4").Value = 30.90100288391113
$ws.Range("F5").Value = 28.24356484413147
$ws.Range("F6").Value = 26.69353222846985
$ws.Range("F7").Value = 26.29553723335266
$ws.Range("F8").Value = 26.31389260292053
$ws.Range("F9").Value = 26.31670951843262
$ws.Range("F10").Value = 26.22270464897156
$ws.Range("F11").Value = 26.42997121810913
$ws.Range("F12").Value = 26.50774264335632
$ws.Range("F13").Value = 26.34843826293945
$ws.Range("F14").Value = 26.39791345596313
$ws.Range("F15").Value = 26.6412320137024
$ws.Range("F16").Value = 26.5207884311676
$ws.Range("F17").Value = 26.52726316452026
$ws.Range("F18").Value = 26.46880984306335
$ws.Range("F19").Value = 26.51732516288757
$ws.Range("F20").Value = 26.51289939880371
$ws.Range("F21").Value = 26.75679087638855

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 26.71737504005432
$ws.Range("F3").Value = 26.37820482254028
$ws.Range("F4").Value = 26.44888615608216
$ws.Range("F5").Value = 26.38553357124329
$ws.Range("F6").Value = 26.44520020484924
$ws.Range("F7").Value = 26.45350241661072
$ws.Range("F8").Value = 26.4281108379364
$ws.Range("F9").Value = 26.42350506782532
$ws.Range("F10").Value = 26.52731561660766
$ws.Range("F11").Value = 26.74628043174744
$ws.Range("F12").Value = 26.55410480499268
$ws.Range("F13").Value = 26.46813440322876
$ws.Range("F14").Value = 26.52019357681274
$ws.Range("F15").Value = 26.62460660934448
$ws.Range("F16").Value = 26.59368944168091
$ws.Range("F17").Value = 26.72600412368774
$ws.Range("F18").Value = 26.6620180606842
$ws.Range("F19").Value = 26.67883896827698
$ws.Range("F20").Value = 26.62525224685669
$ws.Range("F21").Value = 26.89186120033264

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 26.83164644241333
$ws.Range("F3").Value = 26.60200762748718
$ws.Range("F4").Value = 26.5997109413147
$ws.Range("F5").Value = 26.62772154808044
$ws.Range("F6").Value = 26.75187754631042
$ws.Range("F7").Value = 26.6960756778717
$ws.Range("F8").Value = 26.58644080162048
$ws.Range("F9").Value = 26.51874613761902
$ws.Range("F10").Value = 26.54516506195068
$ws.Range("F11").Value = 26.81162762641907
$ws.Range("F12").Value = 26.67332696914673
$ws.Range("F13").Value = 26.55418825149536
$ws.Range("F14").Value = 26.56037092208862
$ws.Range("F15").Value = 26.50875878334045
$ws.Range("F16").Value = 26.41753149032593
$ws.Range("F17").Value = 26.82300734519958
$ws.Range("F18").Value = 26.60252118110657
$ws.Range("F19").Value = 26.63367104530334
$ws.Range("F20").Value = 26.65050005912781
$ws.Range("F21").Value = 26.902508020401

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 26.88038873672485
$ws.Range("F3").Value = 26.60947751998901
$ws.Range("F4").Value = 26.7687623500824
$ws.Range("F5").Value = 26.69063496589661
$ws.Range("F6").Value = 26.67008781433105
$ws.Range("F7").Value = 26.40389919281006
$ws.Range("F8").Value = 26.52817678451538
$ws.Range("F9").Value = 26.5012469291687
$ws.Range("F10").Value = 26.53488731384277
$ws.Range("F11").Value = 26.84052920341492
$ws.Range("F12").Value = 26.64801478385925
$ws.Range("F13").Value = 26.70173811912537
$ws.Range("F14").Value = 26.55263805389404
$ws.Range("F15").Value = 26.56175518035889
$ws.Range("F16").Value = 26.61720991134644
$ws.Range("F17").Value = 26.48142051696777
$ws.Range("F18").Value = 26.49393177032471
$ws.Range("F19").Value = 26.48260402679444
$ws.Range("F20").Value = 26.63510775566101
$ws.Range("F21").Value = 26.8668487071991

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 26.92673206329346
$ws.Range("F3").Value = 26.6506085395813
$ws.Range("F4").Value = 26.61904692649841
$ws.Range("F5").Value = 26.65190100669861
$ws.Range("F6").Value = 26.62372183799744
$ws.Range("F7").Value = 26.67050194740296
$ws.Range("F8").Value = 26.67687273025513
$ws.Range("F9").Value = 26.6266610622406
$ws.Range("F10").Value = 26.63897323608398
$ws.Range("F11").Value = 26.90861678123474
$ws.Range("F12").Value = 26.5458459854126
$ws.Range("F13").Value = 26.66955590248108
$ws.Range("F14").Value = 26.59504294395447
$ws.Range("F15").Value = 26.71591520309448
$ws.Range("F16").Value = 26.56165313720703
$ws.Range("F17").Value = 26.67468285560608
$ws.Range("F18").Value = 26.70881152153015
$ws.Range("F19").Value = 26.64816379547119
$ws.Range("F20").Value = 26.64202523231506
$ws.Range("F21").Value = 26.78550815582276
